$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One of the "Orweb" rows (row 16) was a duplicate of row 15 and is removed.
# Deleting the whole row shifts OSLRMeshTether/SecureSMS/StoryMaker up by one
# row (old rows 17-19 become 16-18), matching the "adding 3 more projects"
# re-numbering shown in the diff.
$ws.Rows("16").Delete()

# The active selection recorded in the sheet view moves to F22.
$ws.Range("F22").Select() | Out-Null

# Page orientation switches from portrait to landscape.
$ws.PageSetup.Orientation = 2
